$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 824, shifting rows 824:865 down to 825:866
$ws.Rows.Item(824).Insert()

# Populate the newly inserted row 824 with the new entry.
# Force column A to stay plain text (like every other date cell in the
# sheet) instead of being auto-recognised as a date value.
$ws.Cells.Item(824, 1).NumberFormat = "@"
$ws.Cells.Item(824, 1).Value = "2026/02/17"
$ws.Cells.Item(824, 1).Style = "Normal"

$ws.Cells.Item(824, 2).Value = "火"
$ws.Cells.Item(824, 3).Value = 6
$ws.Cells.Item(824, 4).Value = 201
